$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Columns Template")
$ws2 = $wb.Worksheets.Item("Data Rules")

# Rename header cells on the "Columns Template" sheet
$ws1.Range("C1").Value = "Middle_Name"
$ws1.Range("G1").Value = "Hours_of_Participation"

# Remove the trailing "Notes" column (H1) — shift remaining cells left
$ws1.Range("H1").Delete(-4159)

# Make "Columns Template" the active sheet/tab again, with G1 selected
$ws1.Activate()
$ws1.Range("G1").Select()
